$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$val)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $val
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.500.79"
Set-TextValue $ws.Range("E2") "  +6.39%  "

Set-TextValue $ws.Range("D3") "2.543.14"
Set-TextValue $ws.Range("E3") "  +6.83%  "

Set-TextValue $ws.Range("E4") "  -0.11%  "

Set-TextValue $ws.Range("D5") "506.00"
Set-TextValue $ws.Range("E5") "  +5.91%  "

Set-TextValue $ws.Range("D6") "159.38"
Set-TextValue $ws.Range("E6") "  +8.53%  "

Set-TextValue $ws.Range("D7") "0.617"

Set-TextValue $ws.Range("E8") "  -0.77%  "

Set-TextValue $ws.Range("D9") "2.587.30"
Set-TextValue $ws.Range("E9") "  +8.67%  "

Set-TextValue $ws.Range("D10") "6.23"
Set-TextValue $ws.Range("E10") "  +14.17%  "

Set-TextValue $ws.Range("E11") "  +7.43%  "

Set-TextValue $ws.Range("E12") "  +6.09%  "

Set-TextValue $ws.Range("E13") "  +1.83%  "

Set-TextValue $ws.Range("D14") "2.975.58"
Set-TextValue $ws.Range("E14") "  +6.25%  "

Set-TextValue $ws.Range("D15") "59.172.18"
Set-TextValue $ws.Range("E15") "  +5.73%  "

Set-TextValue $ws.Range("D16") "22.01"

Set-TextValue $ws.Range("D17") "0.0000138"
Set-TextValue $ws.Range("E17") "  +5.01%  "

Set-TextValue $ws.Range("D18") "2.573.93"
Set-TextValue $ws.Range("E18") "  +7.61%  "

Set-TextValue $ws.Range("D19") "4.74"
Set-TextValue $ws.Range("E19") "  +4.87%  "

Set-TextValue $ws.Range("D20") "333.63"
Set-TextValue $ws.Range("E20") "  +6.31%  "

Set-TextValue $ws.Range("D21") "10.36"
Set-TextValue $ws.Range("E21") "  +6.79%  "

Set-TextValue $ws.Range("D22") "6.08"
Set-TextValue $ws.Range("E22") "  +7.02%  "

Set-TextValue $ws.Range("E23") "  +0.45%  "

Set-TextValue $ws.Range("D24") "60.09"
Set-TextValue $ws.Range("E24") "  +6.04%  "

Set-TextValue $ws.Range("D25") "0.417"
Set-TextValue $ws.Range("E25") "  +5.80%  "

Set-TextValue $ws.Range("E26") "  +6.41%  "

Set-TextValue $ws.Range("D27") "0.996"
Set-TextValue $ws.Range("E27") "  -0.53%  "

Set-TextValue $ws.Range("D28") "2.641.64"
Set-TextValue $ws.Range("E28") "  +5.85%  "

Set-TextValue $ws.Range("D29") "7.57"
Set-TextValue $ws.Range("E29") "  +4.82%  "

Set-TextValue $ws.Range("D30") "0.0₃0832"
Set-TextValue $ws.Range("E30") "  +8.01%  "

Set-TextValue $ws.Range("E31") "  -0.35%  "

Set-TextValue $ws.Range("D32") "19.50"
Set-TextValue $ws.Range("E32") "  +8.44%  "

Set-TextValue $ws.Range("D33") "155.12"
Set-TextValue $ws.Range("E33") "  +5.10%  "

Set-TextValue $ws.Range("E34") "  +5.75%  "

Set-TextValue $ws.Range("E35") "  +9.30%  "

Set-TextValue $ws.Range("E36") "  +8.59%  "

Set-TextValue $ws.Range("D37") "3.93"
Set-TextValue $ws.Range("E37") "  +9.84%  "

Set-TextValue $ws.Range("E38") "  +3.52%  "

Set-TextValue $ws.Range("D39") "3.75"
Set-TextValue $ws.Range("E39") "  +10.85%  "

Set-TextValue $ws.Range("E40") "  +7.16%  "

Set-TextValue $ws.Range("D41") "290.94"
Set-TextValue $ws.Range("E41") "  +14.80%  "

Set-TextValue $ws.Range("D42") "34.80"
Set-TextValue $ws.Range("E42") "  +4.43%  "

Set-TextValue $ws.Range("B43") "Mantle"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.628"
Set-TextValue $ws.Range("E43") "  +8.07%  "

Set-TextValue $ws.Range("B44") "Stellar"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D44") "0.101"
Set-TextValue $ws.Range("E44") "  +7.39%  "

Set-TextValue $ws.Range("D45") "0.0559"
Set-TextValue $ws.Range("E45") "  +4.23%  "

Set-TextValue $ws.Range("D46") "0.994"
Set-TextValue $ws.Range("E46") "  -0.36%  "

Set-TextValue $ws.Range("E47") "  +7.73%  "

Set-TextValue $ws.Range("D48") "19.31"
Set-TextValue $ws.Range("E48") "  +13.74%  "

Set-TextValue $ws.Range("E49") "  +5.25%  "

Set-TextValue $ws.Range("B50") "WhiteBITCoin"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D50") "10.26"
Set-TextValue $ws.Range("E50") "  +0.33%  "

Set-TextValue $ws.Range("B51") "SuiNetwork"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D51") "0.715"
Set-TextValue $ws.Range("E51") "  +14.01%  "
